# Auto-generated Excel COM edit script
# Applies numeric corrections to several Leve-profit rows across sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) per scheduled-runner update.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 8667
$ws.Range("I51").Value = 5000
$ws.Range("J51").Value = 9400.4
$ws.Range("K51").Value = 5000
$ws.Range("L51").Value = 9400.4
$ws.Range("M51").Value = -4516
$ws.Range("N51").Value = -10368.4

$ws.Range("H70").Value = 1713.138
$ws.Range("I70").Value = 735.5
$ws.Range("J70").Value = 4785.7144
$ws.Range("K70").Value = 2206.5
$ws.Range("L70").Value = 14357.1432
$ws.Range("M70").Value = -1936.5
$ws.Range("N70").Value = -14897.1432

$ws.Range("H73").Value = 1713.138
$ws.Range("I73").Value = 735.5
$ws.Range("J73").Value = 4785.7144
$ws.Range("K73").Value = 2206.5
$ws.Range("L73").Value = 14357.1432
$ws.Range("M73").Value = -1270.5
$ws.Range("N73").Value = -16229.1432

$ws.Range("H116").Value = 4220.952
$ws.Range("I116").Value = 3412.7273
$ws.Range("J116").Value = 5110
$ws.Range("K116").Value = 3412.7273
$ws.Range("L116").Value = 5110
$ws.Range("M116").Value = 29.27269999999999
$ws.Range("N116").Value = -11994

$ws.Range("H138").Value = 6325.946
$ws.Range("I138").Value = 4454.1
$ws.Range("J138").Value = 6618.422
$ws.Range("K138").Value = 13362.3
$ws.Range("L138").Value = 19855.266
$ws.Range("M138").Value = -8222.300000000001
$ws.Range("N138").Value = -30135.266

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 935.7143
$ws.Range("I2").Value = 1030.8182
$ws.Range("J2").Value = 587
$ws.Range("K2").Value = 1030.8182
$ws.Range("L2").Value = 587
$ws.Range("M2").Value = -917.8181999999999
$ws.Range("N2").Value = -813

$ws.Range("H116").Value = 935.7143
$ws.Range("I116").Value = 1030.8182
$ws.Range("J116").Value = 587
$ws.Range("K116").Value = 1030.8182
$ws.Range("L116").Value = 587
$ws.Range("M116").Value = 1263.1818
$ws.Range("N116").Value = -5175

$ws.Range("H132").Value = 1680.9375
$ws.Range("I132").Value = 1286.875
$ws.Range("K132").Value = 3860.625
$ws.Range("M132").Value = -1330.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 935.7143
$ws.Range("I3").Value = 1030.8182
$ws.Range("J3").Value = 587
$ws.Range("K3").Value = 1030.8182
$ws.Range("L3").Value = 587
$ws.Range("M3").Value = -916.8181999999999
$ws.Range("N3").Value = -815

$ws.Range("H22").Value = 278.92856
$ws.Range("I22").Value = 219.75
$ws.Range("J22").Value = 357.83334
$ws.Range("K22").Value = 219.75
$ws.Range("L22").Value = 357.83334
$ws.Range("M22").Value = -46.75
$ws.Range("N22").Value = -703.83334

$ws.Range("H38").Value = 33018
$ws.Range("J38").Value = 33018
$ws.Range("L38").Value = 33018
$ws.Range("N38").Value = -33850

$ws.Range("H39").Value = 13526.5
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 13526.5
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 13526.5
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -14304.5

$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()

$ws.Range("H134").Value = 2871.4285
$ws.Range("I134").Value = 2050
$ws.Range("K134").Value = 6150
$ws.Range("M134").Value = -3615

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5908.0557
$ws.Range("I31").Value = 2377.4167
$ws.Range("J31").Value = 12969.333
$ws.Range("K31").Value = 2377.4167
$ws.Range("L31").Value = 12969.333
$ws.Range("M31").Value = -2082.4167
$ws.Range("N31").Value = -13559.333

$ws.Range("H34").Value = 5908.0557
$ws.Range("I34").Value = 2377.4167
$ws.Range("J34").Value = 12969.333
$ws.Range("K34").Value = 2377.4167
$ws.Range("L34").Value = 12969.333
$ws.Range("M34").Value = -2175.4167
$ws.Range("N34").Value = -13373.333

$ws.Range("H58").Value = 2501.138
$ws.Range("I58").Value = 1563.3684
$ws.Range("J58").Value = 4282.9
$ws.Range("K58").Value = 1563.3684
$ws.Range("L58").Value = 4282.9
$ws.Range("M58").Value = -1360.3684
$ws.Range("N58").Value = -4688.9

$ws.Range("H134").Value = 2352.8572
$ws.Range("I134").Value = 2617.5715
$ws.Range("J134").Value = 1294
$ws.Range("K134").Value = 7852.7145
$ws.Range("L134").Value = 3882
$ws.Range("M134").Value = -5317.7145
$ws.Range("N134").Value = -8952

$ws.Range("H136").Value = 2501.138
$ws.Range("I136").Value = 1563.3684
$ws.Range("J136").Value = 4282.9
$ws.Range("K136").Value = 4690.1052
$ws.Range("L136").Value = 12848.7
$ws.Range("M136").Value = -2140.1052
$ws.Range("N136").Value = -17948.7

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1004.41
$ws.Range("J131").Value = 1022.30206
$ws.Range("L131").Value = 3066.90618
$ws.Range("N131").Value = -13146.90618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 15000
$ws.Range("J63").Value = 15000
$ws.Range("L63").Value = 15000
$ws.Range("N63").Value = -16372

$ws.Range("H66").Value = 15000
$ws.Range("J66").Value = 15000
$ws.Range("L66").Value = 45000
$ws.Range("N66").Value = -51864

$ws.Range("H126").Value = 2160.889
$ws.Range("I126").Value = 1814.4615
$ws.Range("J126").Value = 2482.5715
$ws.Range("K126").Value = 5443.3845
$ws.Range("L126").Value = 7447.7145
$ws.Range("M126").Value = -2973.3845
$ws.Range("N126").Value = -12387.7145

$ws.Range("H135").Value = 41250
$ws.Range("J135").Value = 41250
$ws.Range("L135").Value = 41250
$ws.Range("N135").Value = -51390

$ws.Range("H138").Value = 48919.93
$ws.Range("J138").Value = 48919.93
$ws.Range("L138").Value = 48919.93
$ws.Range("N138").Value = -59199.93

$ws.Range("H140").Value = 39744
$ws.Range("J140").Value = 39744
$ws.Range("L140").Value = 39744
$ws.Range("N140").Value = -50104

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 575
$ws.Range("I55").Value = 500
$ws.Range("K55").Value = 500
$ws.Range("M55").Value = -327

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 33700
$ws.Range("J54").Value = 33700
$ws.Range("L54").Value = 33700
$ws.Range("N54").Value = -34740
